# Update the PCRSource sheet header row to match the new PCR model:
# old: forward_primer, reverse_primer, circular, assembly, input, output, type, output_name, id
# new: circular, assembly, input, output, type, output_name, id

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PCRSource")

# Clear the whole header row first (removes forward_primer/reverse_primer and
# shrinks the used range away from columns H and I).
$ws.Rows.Item(1).Clear()

# Write the new header values.
$headers = @("circular", "assembly", "input", "output", "type", "output_name", "id")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
